# Append new OHLCV candle rows (1085-1127) to Sheet1.
# The dataset is a time series of open/high/low/close/volume values; this
# extends it with 43 additional 4h candles (continuing from row 1084).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    ,@(45534.5,0.4209,0.4254,0.399,0.4016,29997851.2)
    ,@(45534.66666666666,0.4017,0.4213,0.4001,0.4157,20533705.9)
    ,@(45534.83333333334,0.4157,0.4219,0.4155,0.421,6691947.5)
    ,@(45535,0.421,0.4256,0.4198,0.4242,6566343.9)
    ,@(45535.16666666666,0.4242,0.427,0.4217,0.4218,7189659.9)
    ,@(45535.33333333334,0.4218,0.4244,0.4181,0.4224,7732274.8)
    ,@(45535.5,0.4224,0.4225,0.4181,0.4214,8106743.8)
    ,@(45535.66666666666,0.4214,0.4233,0.4173,0.4191,9010990.6)
    ,@(45535.83333333334,0.4191,0.4198,0.4179,0.4194,1103924.6)
    ,@(45536,0.4196,0.4203,0.4106,0.4148,5573484.6)
    ,@(45536.16666666666,0.4148,0.4149,0.4057,0.4106,10564896.5)
    ,@(45536.33333333334,0.4107,0.4129,0.4077,0.4127,8079995.5)
    ,@(45536.5,0.4126,0.4135,0.4014,0.4108,12415041.6)
    ,@(45536.66666666666,0.4108,0.4124,0.4072,0.4084,2923714.5)
    ,@(45536.83333333334,0.4132,0.4162,0.3967,0.4023,11692533.2)
    ,@(45537,0.4022,0.409,0.4012,0.4072,7675072.3)
    ,@(45537.16666666666,0.4073,0.4083,0.4002,0.4015,12432196)
    ,@(45537.33333333334,0.4015,0.4076,0.3987,0.4075,20689223.7)
    ,@(45537.5,0.4056,0.4061,0.4005,0.4009,3817458.6)
    ,@(45537.66666666666,0.4033,0.4091,0.4016,0.4066,7045033.6)
    ,@(45537.83333333334,0.4067,0.415,0.4058,0.4108,7958755.5)
    ,@(45538,0.4108,0.415,0.4075,0.4104,11459858.4)
    ,@(45538.16666666666,0.4103,0.4129,0.4094,0.4126,11495386.5)
    ,@(45538.33333333334,0.4127,0.4128,0.4057,0.4079,11392803)
    ,@(45538.5,0.408,0.4127,0.3998,0.4058,23738633.9)
    ,@(45538.66666666666,0.4058,0.4072,0.3982,0.4031,10818004.8)
    ,@(45538.83333333334,0.4032,0.4059,0.3965,0.3978,11709498.1)
    ,@(45539,0.3978,0.4003,0.3781,0.3893,26220407)
    ,@(45539.16666666666,0.3893,0.393,0.3839,0.3896,14801008.5)
    ,@(45539.33333333334,0.3895,0.3895,0.35,0.3738,63738551.2)
    ,@(45539.5,0.3738,0.385,0.3714,0.3845,14406219.7)
    ,@(45539.66666666666,0.3847,0.3921,0.3799,0.3816,12677665.9)
    ,@(45539.83333333334,0.3816,0.3846,0.3791,0.3814,5674776.7)
    ,@(45540,0.3814,0.3879,0.3759,0.3763,8645043.699999999)
    ,@(45540.16666666666,0.3762,0.3788,0.3723,0.3766,5673162.9)
    ,@(45540.33333333334,0.3766,0.3854,0.3637,0.3664,13983117.5)
    ,@(45540.5,0.3664,0.3763,0.362,0.3691,13998396.9)
    ,@(45540.66666666666,0.3691,0.3709,0.3662,0.3675,4665462.2)
    ,@(45540.83333333334,0.3675,0.3721,0.3665,0.3685,4382183.4)
    ,@(45541,0.3686,0.374,0.3625,0.3703,6610582.1)
    ,@(45541.16666666666,0.3702,0.3747,0.3696,0.3701,7253159.5)
    ,@(45541.33333333334,0.37,0.3725,0.3693,0.3705,5138032.4)
    ,@(45541.5,0.3706,0.3707,0.3706,0.3707,4906)
)

$startRow = 1085
$lastRow = $startRow + $newData.Count - 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Column A uses a date/time display style on every existing data row (the
# timestamp column). Copy that formatting from the last pre-existing row
# onto the newly written timestamp cells so the style index matches.
$ws.Range("A1084").Copy() | Out-Null
$ws.Range("A$startRow`:A$lastRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
